$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.821.15"
$ws.Range("E2").Value = "  -1.85%  "

$ws.Range("D3").Value = "1.984.86"
$ws.Range("E3").Value = "  -2.98%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.86"
$ws.Range("E5").Value = "  -0.15%  "

$ws.Range("E6").Value = "  -2.99%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.77"
$ws.Range("E7").Value = "  +9.90%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "58.08"
$ws.Range("E9").Value = "  -4.15%  "

$ws.Range("E10").Value = "  +0.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0739"
$ws.Range("E11").Value = "  -1.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.104"
$ws.Range("E12").Value = "  -2.24%  "

$ws.Range("E13").Value = "  -0.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.66"
$ws.Range("E14").Value = "  -0.09%  "

$ws.Range("D15").Value = "2.276.81"
$ws.Range("E15").Value = "  -2.99%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.32"
$ws.Range("E16").Value = "  -1.90%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.07"
$ws.Range("E17").Value = "  +11.95%  "

$ws.Range("D18").Value = "1.993.49"
$ws.Range("E18").Value = "  -2.70%  "

$ws.Range("D19").Value = "35.745.12"
$ws.Range("E19").Value = "  -1.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.52"
$ws.Range("E20").Value = "  -0.18%  "

$ws.Range("E21").Value = "  -0.82%  "

$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.52"
$ws.Range("E23").Value = "  -1.97%  "

$ws.Range("E24").Value = "  -0.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.62"
$ws.Range("E25").Value = "  +16.62%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.28"
$ws.Range("E26").Value = "  -4.37%  "

$ws.Range("E27").Value = "  +4.93%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "164.74"
$ws.Range("E28").Value = "  -0.19%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.34"
$ws.Range("E29").Value = "  -3.09%  "

$ws.Range("E30").Value = "  -1.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.96"
$ws.Range("E31").Value = "  -1.52%  "

$ws.Range("E32").Value = "  -5.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0975"
$ws.Range("E33").Value = "  +13.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0602"
$ws.Range("E34").Value = "  +1.96%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.41"
$ws.Range("E35").Value = "  -0.80%  "

$ws.Range("E36").Value = "  +10.14%  "

$ws.Range("E37").Value = "  -0.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.78"
$ws.Range("E38").Value = "  -2.34%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.67"
$ws.Range("E39").Value = "  +12.92%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.22"
$ws.Range("E40").Value = "  -0.74%  "

$ws.Range("E41").Value = "  -0.29%  "

$ws.Range("E42").Value = "  -0.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0927"
$ws.Range("E43").Value = "  +2.48%  "

$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("E45").Value = "  +4.10%  "

$ws.Range("E46").Value = "  -0.82%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.72"
$ws.Range("E47").Value = "  +4.60%  "

$ws.Range("D48").Value = "1.361.06"
$ws.Range("E48").Value = "  -3.23%  "

$ws.Range("E49").Value = "  -0.85%  "

$ws.Range("E50").Value = "  +2.42%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.63"
$ws.Range("E51").Value = "  +3.09%  "
